# Updates the cryptos list (Price / Volume(1h) columns, and a couple of
# row re-orderings caused by rank changes) to match the latest scrape,
# as produced by the "Updated cryptos list ... with GitHub Actions" run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value to a cell, forcing it to be stored as text when the
# value looks like a plain number (e.g. "214.07"), so Excel does not silently
# convert it to a numeric cell - matching the original sheet where every
# Price/Volume cell is stored as a string.
function Set-TextValue($range, [string]$value) {
    if ($value -match '^-?\d+(\.\d+)?$') {
        $ws.Range($range).Value = "'" + $value
    } else {
        $ws.Range($range).Value = $value
    }
}

# --- Row 2 (Bitcoin) ---
Set-TextValue "D2" "29.885.53"
$ws.Range("E2").Value = "  +1.22%  "

# --- Row 3 (Ethereum) ---
Set-TextValue "D3" "1.626.10"
$ws.Range("E3").Value = "  +1.90%  "

# --- Row 4 (TetherUSD) ---
$ws.Range("E4").Value = "  -0.11%  "

# --- Row 5 (BNB) ---
Set-TextValue "D5" "214.07"
$ws.Range("E5").Value = "  +1.03%  "

# --- Row 6 (XRP) ---
$ws.Range("E6").Value = "  +1.09%  "

# --- Row 7 (USDC) ---
$ws.Range("E7").Value = "  -0.12%  "

# --- Row 8 (Solana) ---
Set-TextValue "D8" "29.85"
$ws.Range("E8").Value = "  +11.47%  "

# --- Row 10 (Dogecoin) ---
$ws.Range("E10").Value = "  +2.58%  "

# --- Row 11 (TRON) ---
Set-TextValue "D11" "0.0916"
$ws.Range("E11").Value = "  +0.89%  "

# --- Row 12 (WrappedliquidstakedEther2.0) ---
Set-TextValue "D12" "1.859.11"

# --- Row 13 (WrappedEther) ---
Set-TextValue "D13" "1.627.49"
$ws.Range("E13").Value = "  +1.99%  "

# --- Row 14 (Polygon) ---
$ws.Range("E14").Value = "  +6.57%  "

# --- Rows 15-17: ranking shuffled (Chainlink moved up to rank 15,
#     Polkadot down to 16, WrappedBTC down to 17); rank numbers in
#     column A stay the same, only Coin/Link/Price/Volume move.
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D15" "9.19"
$ws.Range("E15").Value = "  +21.84%  "

$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D16" "3.87"
$ws.Range("E16").Value = "  +4.09%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue "D17" "29.898.65"
$ws.Range("E17").Value = "  +1.21%  "

# --- Row 18 (Litecoin) ---
Set-TextValue "D18" "64.95"
$ws.Range("E18").Value = "  +1.97%  "

# --- Row 19 (BitcoinCash) ---
Set-TextValue "D19" "247.96"
$ws.Range("E19").Value = "  +3.08%  "

# --- Row 20 (ShibaInu) ---
$ws.Range("E20").Value = "  +2.02%  "

# --- Row 21 (Dai) ---
$ws.Range("E21").Value = "  -0.06%  "

# --- Row 22 (Uniswap) ---
Set-TextValue "D22" "4.14"
$ws.Range("E22").Value = "  +4.33%  "

# --- Row 23 (Avalanche) ---
Set-TextValue "D23" "9.59"
$ws.Range("E23").Value = "  +4.26%  "

# --- Row 24 (Toncoin) ---
$ws.Range("E24").Value = "  +0.61%  "

# --- Row 25 (Monero) ---
Set-TextValue "D25" "158.62"
$ws.Range("E25").Value = "  +2.39%  "

# --- Row 26 (EthereumClassic) ---
$ws.Range("E26").Value = "  +2.53%  "

# --- Row 27 (Stellar) ---
$ws.Range("E27").Value = "  +2.49%  "

# --- Row 28 (Cosmos) ---
Set-TextValue "D28" "6.59"
$ws.Range("E28").Value = "  +3.37%  "

# --- Row 29 (BinanceUSD) ---
$ws.Range("E29").Value = "  -0.14%  "

# --- Row 30 (Hedera) ---
$ws.Range("E30").Value = "  +2.98%  "

# --- Row 31 (PancakeSwap) ---
$ws.Range("E31").Value = "  +6.07%  "

# --- Row 32 (Filecoin) ---
Set-TextValue "D32" "3.36"
$ws.Range("E32").Value = "  +4.69%  "

# --- Row 33 (InternetComputer(DFINITY)) ---
$ws.Range("E33").Value = "  +1.93%  "

# --- Row 34 (Maker) ---
Set-TextValue "D34" "1.428.00"
$ws.Range("E34").Value = "  -0.17%  "

# --- Row 35 (LidoDAOToken) ---
$ws.Range("E35").Value = "  +7.18%  "

# --- Row 36 (TrustWalletToken) ---
$ws.Range("E36").Value = "  +1.24%  "

# --- Row 37 (MXToken) ---
Set-TextValue "D37" "2.87"
$ws.Range("E37").Value = "  +1.60%  "

# --- Row 38 (HuobiToken) ---
$ws.Range("E38").Value = "  -0.53%  "

# --- Row 39 (VeChain) ---
$ws.Range("E39").Value = "  +3.20%  "

# --- Row 40 (ImmutableX) ---
$ws.Range("E40").Value = "  +3.11%  "

# --- Row 41 (ARBITRUM) ---
Set-TextValue "D41" "0.829"
$ws.Range("E41").Value = "  +3.93%  "

# --- Row 42 (BitcoinSV) ---
Set-TextValue "D42" "55.23"
$ws.Range("E42").Value = "  +3.55%  "

# --- Row 43 (WEMIXToken) ---
$ws.Range("E43").Value = "  +7.83%  "

# --- Rows 44-45: ranking shuffled (RenderToken moved up to rank 44,
#     Kaspa down to 45); rank numbers in column A stay the same.
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D44" "1.98"
$ws.Range("E44").Value = "  +1.17%  "

$ws.Range("B45").Value = "Kaspa"
$ws.Range("C45").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D45" "0.0495"
$ws.Range("E45").Value = "  +0.91%  "

# --- Row 46 (Aave) ---
Set-TextValue "D46" "70.83"
$ws.Range("E46").Value = "  +8.36%  "

# --- Row 47 (PaxDollar) ---
$ws.Range("E47").Value = "  -0.09%  "

# --- Row 48 (FraxShare) ---
Set-TextValue "D48" "5.45"
$ws.Range("E48").Value = "  +2.54%  "

# --- Row 49 (RocketPoolETH) ---
Set-TextValue "D49" "1.767.20"
$ws.Range("E49").Value = "  +1.59%  "

# --- Row 50 (Quant) ---
Set-TextValue "D50" "89.70"
$ws.Range("E50").Value = "  +4.37%  "

# --- Row 51 (BabyDogeCoin) ---
$ws.Range("E51").Value = "  +3.38%  "
